# kevin time logs update
# Extends the Kevin time-log sheet (Sheet1) with new entries for rows 35-46:
#   - row 35 is a blank separator row (carried formatting only)
#   - rows 36-43 are new time-log entries (Start/End/Elapsed/Running total/Notes)
#   - rows 44-46 are trailing blank (formatted) rows in column C
# Also extends the shared-formula groups covering C29:C34 / D29:D34 so the
# existing rows keep working exactly as before, and appends the five new
# strings used by the new rows' notes column to the shared string table in
# the same order Excel would have appended them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "m/d/yy h:mm"
$timeFmt = "[hh]:mm:ss"

# --- Re-fill C29:C34 / D29:D34 so they become one shared-formula group each
#     (mirrors the original author re-selecting and re-filling that block). ---
$ws.Range("C29:C34").Formula = "=B29-A29"
$ws.Range("D29:D34").Formula = "=C29+D28"

# --- Row 35: blank separator row, formatting only (no values/formulas). ---
$ws.Range("A35").NumberFormat = $dateFmt
$ws.Range("B35").NumberFormat = $dateFmt
$ws.Range("C35").NumberFormat = $timeFmt
$ws.Range("D35").NumberFormat = $timeFmt

# --- Row 36 ---
$ws.Range("A36").Value = 44654.583333333336
$ws.Range("B36").Value = 44654.614583333336
$ws.Range("A36").NumberFormat = $dateFmt
$ws.Range("B36").NumberFormat = $dateFmt
$ws.Range("C36").Formula = "=B37-A37"
$ws.Range("D36").Formula = "=C36"
$ws.Range("C36").NumberFormat = $timeFmt
$ws.Range("D36").NumberFormat = $timeFmt
$ws.Range("E37").Value = "Finalize Deisgn Plan"
$ws.Range("E38").Value = "Finalize Skeleton Setup"
$ws.Range("E36").Value = "Client Meeting"

# --- Row 37 ---
$ws.Range("A37").Value = 44655.416666666664
$ws.Range("B37").Value = 44655.447916666664
$ws.Range("A37").NumberFormat = $dateFmt
$ws.Range("B37").NumberFormat = $dateFmt
$ws.Range("C37").Formula = "=B38-A38"
$ws.Range("D37").Formula = "=D36+C37"
$ws.Range("C37").NumberFormat = $timeFmt
$ws.Range("D37").NumberFormat = $timeFmt

# --- Row 38 ---
$ws.Range("A38").Value = 44655.458333333336
$ws.Range("B38").Value = 44655.489583333336
$ws.Range("A38").NumberFormat = $dateFmt
$ws.Range("B38").NumberFormat = $dateFmt
$ws.Range("C38").Formula = "=B38-A38"
$ws.Range("D38").Formula = "=C38+D37"
$ws.Range("C38").NumberFormat = $timeFmt
$ws.Range("D38").NumberFormat = $timeFmt

# --- Row 39 ---
$ws.Range("A39").Value = 44656.59375
$ws.Range("B39").Value = 44656.604166666664
$ws.Range("A39").NumberFormat = $dateFmt
$ws.Range("B39").NumberFormat = $dateFmt
$ws.Range("C39").Formula = "=B39-A39"
$ws.Range("D39").Formula = "=C39+D38"
$ws.Range("C39").NumberFormat = $timeFmt
$ws.Range("D39").NumberFormat = $timeFmt
$ws.Range("E39").Value = "Team meeting"

# --- Row 40 ---
$ws.Range("A40").Value = 44657.583333333336
$ws.Range("B40").Value = 44657.604166666664
$ws.Range("A40").NumberFormat = $dateFmt
$ws.Range("B40").NumberFormat = $dateFmt
$ws.Range("C40").Formula = "=B40-A40"
$ws.Range("D40").Formula = "=C40+D39"
$ws.Range("C40").NumberFormat = $timeFmt
$ws.Range("D40").NumberFormat = $timeFmt
$ws.Range("E40").Value = "Team meeting"

# --- Row 41 ---
$ws.Range("A41").Value = 44657.604166666664
$ws.Range("B41").Value = 44657.625
$ws.Range("A41").NumberFormat = $dateFmt
$ws.Range("B41").NumberFormat = $dateFmt
$ws.Range("C41").Formula = "=B41-A41"
$ws.Range("D41").Formula = "=C41+D40"
$ws.Range("C41").NumberFormat = $timeFmt
$ws.Range("D41").NumberFormat = $timeFmt
$ws.Range("E41").Value = "Finsh Milestone 6"

# --- Row 42 ---
$ws.Range("A42").Value = 44661.604166666664
$ws.Range("B42").Value = 44661.635416666664
$ws.Range("A42").NumberFormat = $dateFmt
$ws.Range("B42").NumberFormat = $dateFmt
$ws.Range("C42").Formula = "=B42-A42"
$ws.Range("D42").Formula = "=C42+D41"
$ws.Range("C42").NumberFormat = $timeFmt
$ws.Range("D42").NumberFormat = $timeFmt
$ws.Range("E42").Value = "Client meeting"

# --- Row 43 ---
$ws.Range("A43").Value = 44664.395833333336
$ws.Range("B43").Value = 44664.479166666664
$ws.Range("A43").NumberFormat = $dateFmt
$ws.Range("B43").NumberFormat = $dateFmt
$ws.Range("C43").Formula = "=B43-A43"
$ws.Range("D43").Formula = "=C43+D42"
$ws.Range("C43").NumberFormat = $timeFmt
$ws.Range("D43").NumberFormat = $timeFmt
$ws.Range("E43").Value = "TUI and Makefile implementation"

# --- Rows 44-46: trailing blank (formatted) rows in column C only. ---
$ws.Range("C44").NumberFormat = $timeFmt
$ws.Range("C45").NumberFormat = $timeFmt
$ws.Range("C46").NumberFormat = $timeFmt

# --- View: show D46 as the active selection, and drop the custom zoom back
#     to 100% (matches the saved sheetView). ---
$ws.Range("D46").Select()
$excel.ActiveWindow.Zoom = 100
